# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.466.67'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '1.877.26'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("E4").Value = '  +0.54%  '
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.20'
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  +3.42%  '
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = '  +1.00%  '
$ws.Range("E7").Value = '  +0.45%  '
$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.81'
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = '  +7.52%  '
$ws.Range("E9").Value = '  +0.30%  '
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0699'
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("E12").Value = '  +1.88%  '
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.61'
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = '1.881.77'
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("E16").Value = '  +2.20%  '
$ws.Range("D17").Value = '35.436.93'
$ws.Range("E17").Value = '  +1.39%  '
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '71.34'
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").Value = '0.0₃0803'
$ws.Range("E19").Value = '  +2.09%  '
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '243.51'
$ws.Range("D20").Style = $__style
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.36'
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = '  +1.39%  '
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.76'
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("E24").Value = '  -0.67%  '
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.45'
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.90'
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = '  +23.85%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.23'
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = '  +5.30%  '
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.83'
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("E29").Value = '  +1.46%  '
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0564'
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = '  +2.13%  '
$ws.Range("E31").Value = '  +2.46%  '
$ws.Range("E33").Value = '  +2.28%  '
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.80'
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = '  +21.69%  '
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.838'
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = '  +19.38%  '
$ws.Range("E36").Value = '  +6.32%  '
$ws.Range("E37").Value = '  +7.31%  '
$ws.Range("E38").Value = '  +3.40%  '
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0205'
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = '  +5.13%  '
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '91.05'
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").Value = '1.353.33'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.28'
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = '  +3.16%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0601'
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = '  +14.75%  '
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.35'
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = '  +2.59%  '
$ws.Range("E45").Value = '  +55.11%  '
$ws.Range("E46").Value = '  +0.46%  '
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.66'
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = '  +6.58%  '
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.72'
$ws.Range("D48").Style = $__style
$ws.Range("D49").Value = '2.059.96'
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("E50").Value = '  +3.29%  '
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.44'
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = '  +0.62%  '
